# cli progress, i guess
#
# - Rename the lone worksheet from "names" to the generic "Sheet1"
# - Select column B (whole column) so the persisted view shows
#   activeCell="B1" sqref="B1:B1048576"
# - Strip every built-in/custom cell style except "Normal" so the
#   style table collapses down to the bare minimum

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Sheet1"

$styleNames = @()
foreach ($s in $wb.Styles) {
    $styleNames += $s.Name
}
foreach ($n in $styleNames) {
    if ($n -ne "Normal") {
        $wb.Styles.Item($n).Delete()
    }
}

$ws.Range("B:B").Select()
